$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "class:scoreExample" -> "class:satisfactionImportance" in the first
#    content paragraph, kept as the separate runs "class:" and
#    "satisfactionImportance" (matching the surrounding "TTT" runs).
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("scoreExample") | Out-Null
$rng.Text = "satisfactionImportance"

# Use transient bookmarks to force a run boundary at the "class:" /
# "satisfactionImportance" seam (and at the end, before the trailing "TTT")
# without leaving any bookmark or formatting residue behind once removed.
$splitBefore = $d.Range($rng.Start, $rng.Start)
$d.Bookmarks.Add("TempSplitA", $splitBefore) | Out-Null
$splitAfter = $d.Range($rng.End, $rng.End)
$d.Bookmarks.Add("TempSplitB", $splitAfter) | Out-Null
$d.Bookmarks("TempSplitA").Delete()
$d.Bookmarks("TempSplitB").Delete()

# ---------------------------------------------------------------------------
# 2) Relocate the "_GoBack" bookmark: drop it from the first paragraph and
#    drop the "class:percentageExample" text from the paragraph that used to
#    hold it, leaving that paragraph empty.
# ---------------------------------------------------------------------------
$d.Bookmarks("_GoBack").Delete()

$pctPara = $d.Paragraphs(5).Range
$pctPara.MoveEnd(1, -1) | Out-Null
$pctPara.Delete()

# ---------------------------------------------------------------------------
# 3) Remove the "class:reportmark" paragraph's text too (its own paragraph
#    mark, and the blank paragraphs around it, are merged away below).
# ---------------------------------------------------------------------------
$reportPara = $d.Paragraphs(7).Range
$reportPara.MoveEnd(1, -1) | Out-Null
$reportPara.Delete()

# Merge the (now all-empty) paragraph marks for paragraphs 5-8 forward so
# only a single empty paragraph remains in their place, immediately before
# the section break.
$d.Paragraphs(5).Range.Delete()
$d.Paragraphs(5).Range.Delete()
$d.Paragraphs(5).Range.Delete()
$d.Paragraphs(5).Range.Delete()

# ---------------------------------------------------------------------------
# 4) Re-add the "_GoBack" bookmark, now collapsed inside that final, empty
#    paragraph.
# ---------------------------------------------------------------------------
$finalRange = $d.Paragraphs(5).Range
$finalRange.MoveEnd(1, -1) | Out-Null
$d.Bookmarks.Add("_GoBack", $finalRange) | Out-Null
